# Applies the "mlk.docx" edit:
#  1. Append two trailing spaces to the first sentence.
#  2. Append a red-colored parenthetical note, inserted as three separate
#     runs (mirrors how the source edit split the text across runs).
#  3. Add a new, empty, shaded (F9F9F9) paragraph at the very end of the
#     document, just before the final section break.

$d = $word.ActiveDocument

# --- Step 1: first paragraph text -----------------------------------------

# 1a. Append two trailing spaces to the original sentence.
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "This is a Microsoft word document.  ", 2) | Out-Null

# 1b. Insert three red-colored runs right after it (split into 3 runs, as
#     in the source edit: "(This is a change – Ve" | "rsion for main branch" | ")")
$p1 = $d.Paragraphs(1)
$end = $p1.Range
$end.Collapse(0)     | Out-Null  # wdCollapseEnd
$end.MoveEnd(1, -1)  | Out-Null  # step back before the paragraph mark (wdCharacter)
$end.InsertAfter("(This is a change " + [char]0x2013 + " Ve")
$end.Font.Color = 255   # wdColorRed -> w:val="FF0000"

$p1 = $d.Paragraphs(1)
$end = $p1.Range
$end.Collapse(0)    | Out-Null
$end.MoveEnd(1, -1) | Out-Null
$end.InsertAfter("rsion for main branch")
$end.Font.Color = 255

$p1 = $d.Paragraphs(1)
$end = $p1.Range
$end.Collapse(0)    | Out-Null
$end.MoveEnd(1, -1) | Out-Null
$end.InsertAfter(")")
$end.Font.Color = 255

# --- Step 2: new blank shaded paragraph at the very end of the document ---
$endOfDoc = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endOfDoc.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:shd w:val='clear' w:color='auto' w:fill='F9F9F9'/></w:pPr></w:p>") | Out-Null

Write-Output "Done."
